# Automatische test-sync: 2025-06-26 19:33:50
# Append the new log row (row 11) to the "Logs" sheet and refresh the
# "Dashboard" sheet's category count for "Bestelling / Levering".

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append row 11 -------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(11, 1).Value2 = "Kun je 2 dozen nitrilhandschoenen bestellen?"
$logs.Cells.Item(11, 2).Value2 = "MailMind Test <mailmind.test@zohomail.eu>"
$logs.Cells.Item(11, 3).Value2 = "He Johan,`nKun je 2 dozen nitrilhandschoenen bestellen?`nMarc`nSent using {0}"
$logs.Cells.Item(11, 4).Value2 = "Bestelling / Levering"
$logs.Cells.Item(11, 5).Value2 = "Beste Marc,`nBedankt voor je e-mail. Om je verzoek voor het bestellen van 2 dozen nitrilhandschoenen te verwerken, heb ik wat meer informatie nodig. Kun je mij laten weten welke maat(en) nitrilhandschoenen je nodig hebt? Zodra ik deze informatie van je heb ontvangen, zal ik de bestelling voor je plaatsen.`nMet vriendelijke groet,`nJohan  `nE-mailassistent"
$logs.Cells.Item(11, 6).Value2 = "2025-06-26 19:33:20"
$logs.Cells.Item(11, 7).Value2 = "Ja"
$logs.Cells.Item(11, 8).Value2 = "Nee"
$logs.Cells.Item(11, 9).Value2 = "Ja"

# The multi-line content above makes the COM layer auto-size the row;
# restore it to the sheet's default (un-customized) row height so the
# row tag matches the rest of the sheet.
$logs.Rows.Item(11).AutoFit()

# Extend the conditional formatting ranges from row 2:10 to row 2:11 so
# the newly added row keeps getting highlighted like the rest of the log.
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))

# --- Dashboard sheet: bump the "Bestelling / Levering" tally -------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value2 = 6
